# 29-01-2019 changes to add country wise DF
$wb = $excel.ActiveWorkbook

$wsCountry = $wb.Worksheets.Item("Table-5.1")
$wsCountry.Range("C7").Value = "Others"
$wsCountry.Range("D7").Value = "Others"
$wsCountry.Range("E7").Value = "Others"

$wsTable1 = $wb.Worksheets.Item("Table -1.1")
$wsTable1.Range("C5").Value = 66368

$wsCountry.Activate()
$wsCountry.Range("J15").Select()
